$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")

# Swap Runmode values in rows 13 and 15 (column C)
$ws1.Range("C13").Value = "N"
$ws1.Range("C15").Value = "Y"

# Update selection on sheet1
$ws1.Range("E18").Select()
